$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update font name in row 2 and remove now-unused rows 3 and 4 (Roboto, Georgia)
$ws.Range("A2").Value = "Comic Sans"

# Clear the now-unneeded rows (previously Roboto / Georgia entries)
$ws.Range("A3:B4").Delete()

# Update the sheet selection to match the saved view state
$ws.Range("J5:K5").Select()
